# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
#  * New "Player Info" sheet inserted as the first sheet, ahead of the
#    existing "ODI Batting" / "ODI Bowling" sheets.
#  * On "ODI Batting" and "ODI Bowling", the MATCH_CARD_LINK column is
#    renamed to MATCH_CODE and its values trimmed down from the full
#    scorecard URL to just the bare match-code number.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet ahead of the current first sheet
#    (Worksheets.Add() with no args inserts before the active sheet, and
#    the workbook's active sheet is the current first sheet, "ODI
#    Batting" -- so this lands it exactly where the diff expects it).
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Keep the player ID as text (matches the source sheets, which store
# every value -- numeric-looking or not -- as a string).
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5664"
$playerInfo.Range("B2").Value = "Muhammad Musa Khan"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# Match the bold / centered / top-aligned / thin-bordered header style
# used by the other two sheets' header rows.
$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous
$hdr.Borders.Weight = 2            # xlThin

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, URLs -> codes
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2:D3").NumberFormat = "@"
$batting.Range("D2").Value = "4433"
$batting.Range("D3").Value = "4434"

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, URLs -> codes
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Range("B2:B3").NumberFormat = "@"
$bowling.Range("B2").Value = "4433"
$bowling.Range("B3").Value = "4434"

Write-Output "Player Info sheet added; MATCH_CARD_LINK -> MATCH_CODE done."
